# Update to latest job ad data / data text and sources (3-1_DataTable, Sheet1)
#
# Row 13 holds the "Job adverts by occupation" data source row. Its
# "Latest period (release date)" (C13) and "Next period (release date)" (D13)
# values are refreshed to the newer reporting period.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "Mar 2025 (25/04/25)"
$ws.Range("D13").Value = "Apr 2025 (May 2025)"

# Reflect the author's final view/selection state: scrolled up one row with
# cell D14 selected.
$win = $excel.ActiveWindow
$ws.Range("D14").Select()
$win.ScrollRow = 2
$win.ScrollColumn = 1
